$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 16.45680000000001
$ws.Range("C3").Value = -12.28369999999999
$ws.Range("E3").Value = 16.90489999999999
$ws.Range("C4").Value = -12.28920000000001
$ws.Range("E9").Value = 17.24370000000001
$ws.Range("A11").Value = -21.5305
$ws.Range("A12").Value = -21.5526
$ws.Range("C14").Value = -13.49339999999999
$ws.Range("A15").Value = -21.69360000000001
$ws.Range("E15").Value = 16.1123
$ws.Range("E19").Value = 16.57870000000001
$ws.Range("E20").Value = 15.83559999999999
$ws.Range("E25").Value = 16.89890000000001
$ws.Range("C26").Value = -12.93350000000001
$ws.Range("A27").Value = -22.00759999999999
$ws.Range("E27").Value = 16.7971
$ws.Range("A28").Value = -21.93199999999999
$ws.Range("E28").Value = 16.67290000000001
$ws.Range("E30").Value = 16.01100000000001
$ws.Range("A31").Value = -21.5618
$ws.Range("C31").Value = -13.3969
$ws.Range("A32").Value = -21.30199999999999
$ws.Range("E32").Value = 16.71929999999999
$ws.Range("C35").Value = -13.34230000000001
$ws.Range("A36").Value = -20.7815
$ws.Range("C37").Value = -14.57549999999999
$ws.Range("A38").Value = -19.1437
$ws.Range("C39").Value = -12.2151
$ws.Range("C40").Value = -13.378
$ws.Range("E44").Value = 16.70780000000001
$ws.Range("C45").Value = -13.4944
$ws.Range("A46").Value = -21.49939999999999
$ws.Range("E47").Value = 16.6575
$ws.Range("C52").Value = -10.7409
$ws.Range("A54").Value = -21.70469999999999
$ws.Range("A55").Value = -22.31590000000001
$ws.Range("A56").Value = -22.18130000000001
$ws.Range("C57").Value = -14.4429
$ws.Range("E58").Value = 16.39040000000001
$ws.Range("E62").Value = 16.7056
$ws.Range("A67").Value = -21.46059999999998
$ws.Range("A69").Value = -21.70699999999998
$ws.Range("A72").Value = -21.45379999999999
$ws.Range("A73").Value = -19.91789999999999
$ws.Range("E77").Value = 17.40410000000002
$ws.Range("E78").Value = 16.39980000000002
$ws.Range("C81").Value = -13.0338
$ws.Range("A83").Value = -21.8086
$ws.Range("C83").Value = -12.912
$ws.Range("E84").Value = 16.94469999999999
$ws.Range("A86").Value = -22.2163
$ws.Range("E89").Value = 17.21930000000001
$ws.Range("A91").Value = -21.47630000000001
$ws.Range("E91").Value = 17.92990000000002
$ws.Range("E92").Value = 18.04910000000002
$ws.Range("A93").Value = -21.14429999999999
$ws.Range("E96").Value = 16.76659999999999
$ws.Range("A99").Value = -20.06289999999999
$ws.Range("C100").Value = -13.49839999999999
$ws.Range("C102").Value = -11.78040000000001
$ws.Range("E102").Value = 16.79870000000001
